# Standard User name change in Contacts test data files - 6th Mar 2024
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

# Update the Standard User name value
$ws.Range("A2").Value = "Ayati Arvind"

# Reflect the new selection recorded on save (D8, even though sheet data is only A1:A2)
$ws.Range("D8").Select()
